# WorkingHoursLog.xlsx - "Added menu item for helsinkikanava videos."
#
# Row 12 (2015-08-26) gets a logged entry: 2 hours, "Bug fixes."
# The Total formula in C24 (=SUM(C2:C22)) recalculates automatically.
# The sheet's selection moves to C13 and the view scrolls right one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log the new working-hours entry on row 12.
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = "Bug fixes."

# Update the selected cell and scroll the window so column B is leftmost.
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
